$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '307.50'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '-4.62%'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '39.99'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '-6.40%'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.014'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '-4.47%'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.07669'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '4.234'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '-2.74%'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.618'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '-10.03%'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.8892'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '-6.93%'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.1004'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '-10.92%'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1737'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '-6.62%'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08924'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '-4.93%'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.04387'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '-5.05%'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.1056'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '-0.35%'
$ws.Range("B14").Value = 'TigerCash'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.005826'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '-1.84%'
$ws.Range("B15").Value = 'LEO'
$ws.Range("C15").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.357'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '-0.68%'
$ws.Range("B16").Value = 'BTSEToken'
$ws.Range("C16").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.530'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '0.49%'
$ws.Range("B17").Value = 'BitpandaEcosystemToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.3361'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '-0.08%'
$ws.Range("B18").Value = 'MCDex'
$ws.Range("C18").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.992'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '-6.34%'
$ws.Range("B19").Value = 'ProBitToken'
$ws.Range("C19").Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.1342'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '-1.82%'
$ws.Range("B20").Value = 'ZBToken'
$ws.Range("C20").Value = 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.3030'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '11.72%'
$ws.Range("B21").Value = 'BitForexToken'
$ws.Range("C21").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.001271'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '-0.82%'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.04228'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '0.92%'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.001201'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '-4.37%'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.004066'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '-5.77%'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0001222'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '-6.77%'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '-0.47%'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02354'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '-9.40%'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05166'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '-5.70%'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.007940'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '1.36%'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.1325'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '-4.79%'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.006568'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '-0.51%'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.002003'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '-5.91%'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.007629'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '-12.06%'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.3052'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '-11.49%'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00006580'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '-5.84%'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.00000000752'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '-0.35%'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.003474'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '-0.58%'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '40.93%'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.00002104'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '-0.35%'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0002004'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '-0.35%'
